# Weekly data refresh for "Fruta, Vega Modelo de Temuco - Naranja":
# insert 4 new daily-report rows just above the existing row 693, which
# pushes the previously-existing rows 693-700 down to 697-704.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 693 (shifts old rows 693:700 down to 697:704).
$ws.Range("A693:A696").EntireRow.Insert()

# --- New row 693: Navel / Calibre 88, imported from the US ---
$ws.Cells.Item(693, 1).Value = 10
$ws.Cells.Item(693, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(693, 3).Value = "La Araucanía"
$ws.Cells.Item(693, 4).Value = 44628
$ws.Cells.Item(693, 5).Value = 9
$ws.Cells.Item(693, 6).Value = "Fruta"
$ws.Cells.Item(693, 7).Value = 100102
$ws.Cells.Item(693, 8).Value = "Cítricos"
$ws.Cells.Item(693, 9).Value = 100102005
$ws.Cells.Item(693, 10).Value = "Naranja"
$ws.Cells.Item(693, 11).Value = "Navel"
$ws.Cells.Item(693, 12).Value = "Calibre 88"
$ws.Cells.Item(693, 13).Value = 200
$ws.Cells.Item(693, 14).Value = 28000
$ws.Cells.Item(693, 15).Value = 30000
$ws.Cells.Item(693, 16).Value = 29000
$ws.Cells.Item(693, 17).Value = "`$/caja 18 kilos importada"
$ws.Cells.Item(693, 18).Value = "EE.UU."
$ws.Cells.Item(693, 19).Value = 1611
$ws.Cells.Item(693, 20).Value = 18

# --- New row 694: Valencia / Especial, sold in bins ---
$ws.Cells.Item(694, 1).Value = 10
$ws.Cells.Item(694, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(694, 3).Value = "La Araucanía"
$ws.Cells.Item(694, 4).Value = 44628
$ws.Cells.Item(694, 5).Value = 9
$ws.Cells.Item(694, 6).Value = "Fruta"
$ws.Cells.Item(694, 7).Value = 100102
$ws.Cells.Item(694, 8).Value = "Cítricos"
$ws.Cells.Item(694, 9).Value = 100102005
$ws.Cells.Item(694, 10).Value = "Naranja"
$ws.Cells.Item(694, 11).Value = "Valencia"
$ws.Cells.Item(694, 12).Value = "Especial"
$ws.Cells.Item(694, 13).Value = 2
$ws.Cells.Item(694, 14).Value = 300000
$ws.Cells.Item(694, 15).Value = 300000
$ws.Cells.Item(694, 16).Value = 300000
$ws.Cells.Item(694, 17).Value = "`$/bins (400 kilos)"
$ws.Cells.Item(694, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(694, 19).Value = 750
$ws.Cells.Item(694, 20).Value = 400

# --- New row 695: Valencia / Especial, boxed ---
$ws.Cells.Item(695, 1).Value = 10
$ws.Cells.Item(695, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(695, 3).Value = "La Araucanía"
$ws.Cells.Item(695, 4).Value = 44628
$ws.Cells.Item(695, 5).Value = 9
$ws.Cells.Item(695, 6).Value = "Fruta"
$ws.Cells.Item(695, 7).Value = 100102
$ws.Cells.Item(695, 8).Value = "Cítricos"
$ws.Cells.Item(695, 9).Value = 100102005
$ws.Cells.Item(695, 10).Value = "Naranja"
$ws.Cells.Item(695, 11).Value = "Valencia"
$ws.Cells.Item(695, 12).Value = "Especial"
$ws.Cells.Item(695, 13).Value = 140
$ws.Cells.Item(695, 14).Value = 18000
$ws.Cells.Item(695, 15).Value = 19000
$ws.Cells.Item(695, 16).Value = 18571
$ws.Cells.Item(695, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(695, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(695, 19).Value = 1238
$ws.Cells.Item(695, 20).Value = 15

# --- New row 696: Valencia / Primera, tray ---
$ws.Cells.Item(696, 1).Value = 10
$ws.Cells.Item(696, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(696, 3).Value = "La Araucanía"
$ws.Cells.Item(696, 4).Value = 44628
$ws.Cells.Item(696, 5).Value = 9
$ws.Cells.Item(696, 6).Value = "Fruta"
$ws.Cells.Item(696, 7).Value = 100102
$ws.Cells.Item(696, 8).Value = "Cítricos"
$ws.Cells.Item(696, 9).Value = 100102005
$ws.Cells.Item(696, 10).Value = "Naranja"
$ws.Cells.Item(696, 11).Value = "Valencia"
$ws.Cells.Item(696, 12).Value = "Primera"
$ws.Cells.Item(696, 13).Value = 280
$ws.Cells.Item(696, 14).Value = 15000
$ws.Cells.Item(696, 15).Value = 16000
$ws.Cells.Item(696, 16).Value = 15500
$ws.Cells.Item(696, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(696, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(696, 19).Value = 1033
$ws.Cells.Item(696, 20).Value = 15
